# Split the run containing "БАФ-М ... Р-123 на одну антенну..." into three runs,
# inserting an extra "М" right after "Р-123" (i.e. Р-123 -> Р-123М),
# matching commit "add counter for standarts" (adds the "М" index/counter used
# elsewhere in the doc for the Р-123М radio station designation).

$d = $word.ActiveDocument

$oldText = "- блок антенных фильтров (БАФ—М) для совместной работы двух радиостанций Р-123 на одну антенну с высокочастотным кабелем;"

$target = $null
$count = $d.Paragraphs.Count
for ($i = 1; $i -le $count; $i++) {
    $p = $d.Paragraphs.Item($i)
    $t = $p.Range.Text
    if ($t.TrimEnd([char]13, [char]7) -eq $oldText) {
        $target = $p.Range
        break
    }
}

if ($target -eq $null) {
    # Fallback: match on a unique substring if the exact/full-text compare
    # above didn't line up (e.g. whitespace/line-ending quirks).
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        $t = $p.Range.Text
        if ($t -like "*БАФ*Р-123 на одну антенну*") {
            $target = $p.Range
            break
        }
    }
}

if ($target -eq $null) {
    throw "Could not locate target paragraph for БАФ-М edit"
}

$xml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p w:rsidR="0052309E" w:rsidRPr="006805B9" w:rsidRDefault="0052309E" w:rsidP="0052309E"><w:pPr><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr></w:pPr><w:r w:rsidRPr="006805B9"><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>- блок антенных фильтров (БАФ—М) для совместной работы двух радиостанций Р-123</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t>М</w:t></w:r><w:r><w:rPr><w:rFonts w:ascii="Times New Roman" w:hAnsi="Times New Roman" w:cs="Times New Roman"/><w:sz w:val="26"/><w:szCs w:val="26"/></w:rPr><w:t xml:space="preserve"> на одну антенну с высокочастотным кабелем;</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$target.InsertXML($xml)
